$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates derived from the commit diff (refreshed crypto price/ranking data).
# Cells whose new value looks purely numeric are pre-formatted as Text ("@") so Excel
# stores the exact original string (matching the source inlineStr values) instead of
# silently converting them to floating point numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.43'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.605'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05619'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.472'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8030'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.067'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0005747'
$ws.Range('E10').Value = '9OneONE'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1427'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07463'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03186'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.02985'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09259'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001679'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'MCDex'
$ws.Range('C17').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.250'
$ws.Range('E17').Value = '16MCDexMCB'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04690'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006269'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.001049'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.003811'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001502'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0004605'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.981'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.119'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1277'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04179'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007129'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003504'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1045'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009805'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005629'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6808'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.02871'
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
